# M12 Froze Token Embeddings + Decoder 1
# Update B-column tokens and C-column counts on the "LJ Speech" sheet
# to reflect the latest ASR results run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LJ Speech")

# Row 2
$ws.Range("B2").Value = "<they>"
$ws.Range("C2").Value = 56

# Row 3
$ws.Range("C3").Value = 56

# Row 4
$ws.Range("B4").Value = "<is>"
$ws.Range("C4").Value = 57

# Row 5
$ws.Range("B5").Value = "<is>"
$ws.Range("C5").Value = 50

# Row 6
$ws.Range("C6").Value = 56

# Row 7
$ws.Range("C7").Value = 50

# Row 8
$ws.Range("B8").Value = "<word>"
$ws.Range("C8").Value = 51

# Row 9
$ws.Range("C9").Value = 52

# Row 10
$ws.Range("C10").Value = 53

# Row 11
$ws.Range("C11").Value = 51

# Row 12
$ws.Range("C12").Value = 50

# Row 13
$ws.Range("B13").Value = "<other>"

# Row 14
$ws.Range("C14").Value = 51

# Row 15
$ws.Range("B15").Value = "<my>"
$ws.Range("C15").Value = 54
